$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update aggregate statistics after closing Trade #41
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1197.59   # Current Capital
$summary.Range("B4").Value = -2.41    # Total P&L $
$summary.Range("B5").Value = -1.18    # Total P&L %
$summary.Range("B6").Value = 41       # Total Trades
$summary.Range("B8").Value = 22       # Losing Trades
$summary.Range("B9").Value = 39.02    # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update the MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 97.59     # Capital
$status.Range("D4").Value = 41        # Trades
$status.Range("E4").Value = -2.41     # P&L $
$status.Range("F4").Value = -2.41     # P&L %
$status.Range("G4").Value = 39.02     # Win Rate %

# ---------------------------------------------------------------------------
# Sheets "All Trades" and "MarketMaking": append the newly closed Trade #41
# ---------------------------------------------------------------------------
$newRow = 42

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item($newRow, 1).Value = 41

    # Date/Time columns must stay plain text (matching the rest of the
    # table) instead of being auto-converted into date/time serials.
    $dateCell = $ws.Cells.Item($newRow, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.Style = "Normal"

    $timeCell = $ws.Cells.Item($newRow, 3)
    $timeCell.NumberFormat = "@"
    $timeCell.Value = "13:27:17"
    $timeCell.Style = "Normal"

    $ws.Cells.Item($newRow, 4).Value = "MarketMaking"
    $ws.Cells.Item($newRow, 5).Value = "UP"
    $ws.Cells.Item($newRow, 6).Value = 0.75
    $ws.Cells.Item($newRow, 7).Value = 0.65
    $ws.Cells.Item($newRow, 8).Value = "CLOSED"
    $ws.Cells.Item($newRow, 9).Value = -13.3333
    $ws.Cells.Item($newRow, 10).Value = -0.1
    $ws.Cells.Item($newRow, 11).Value = 97.59
    $ws.Cells.Item($newRow, 12).Value = 0
    $ws.Cells.Item($newRow, 13).Value = 0
    $ws.Cells.Item($newRow, 14).Value = 0.6
    $ws.Cells.Item($newRow, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($newRow, 16).Value = "early_exit"
    $ws.Cells.Item($newRow, 17).Value = 0.13
}
